# feat: add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new row for 2022-Q3 above the
#    existing 2022-Q2 / 2022-Q1 rows (they shift down by one row).
# 2. Duplicate the existing "2022-Q2" sheet (placed right before it,
#    which is exactly where the new tab belongs) and rename the copy
#    "2022-Q3", then update its fund-position figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) "总计" sheet - shift the two existing rows down and insert the
#    new 2022-Q3 row at the top of the data block.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room for the new row 4 by copying row 3's formatting (style)
# down into row 4 first (keeps the bordered/bold "A" column style).
$summary.Range("A3").Copy($summary.Range("A4"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.39

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.66

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.6

# ---------------------------------------------------------------
# 2) Add the "2022-Q3" sheet by duplicating "2022-Q2" (same layout),
#    which Excel places immediately before the source sheet - i.e.
#    right where the new quarter tab belongs - then edit its values.
# ---------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# D2:G2 hold text-formatted figures (e.g. "4.45") in the source data,
# so force the range to Text before writing to avoid Excel silently
# re-typing them as numbers.
$q3.Range("D2:G2").NumberFormat = "@"
$q3.Range("D2").Value = "4.45"
$q3.Range("E2").Value = "96.33"
$q3.Range("F2").Value = "8.79"
$q3.Range("G2").Value = "0.3912"
$q3.Range("H2").Value = 7
